$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the student record values (row 2) ---
$ws.Range("A2").Value = "Dharun Vignesh G"
$ws.Range("B2").Value = 26
$ws.Range("E2").Value = 8
$ws.Range("G2").Value = 1111
$ws.Range("J2").Value = "AAV"

# --- Email cell: new text + refreshed hyperlink (same target, new display/tooltip) ---
$ws.Range("I2").Value = "rfg@gmail"
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:gdvbca@gmail.com", [Type]::Missing, "mailto:rfg@gmail", "rfg@gmail")
$ws.Range("I2").Font.Color = 8388736

# --- B1 header loses its (accidental) Hyperlink look, back to the default style ---
$ws.Range("B1").ClearFormats()

# --- Move the active selection to G2 ---
[void]$ws.Range("G2").Select()
